$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast error values (naive component forecaster bug fix)
$data = @{
    2  = @(0.2109899894693604, 0.6674624635798462, 1.146324606009462, 1.070665496786677, 1.060115081941694, 51)
    3  = @(0.04764615694829157, 0.6468917687613973, 1.115891656663563, 1.056357731388171, 1.065996467999481, 50)
    4  = @(0.180187109103845, 0.6632043662894834, 1.225949235404478, 1.107225918864112, 1.103787138057222, 49)
    5  = @(0.07615678762443598, 0.6620181254576488, 1.170276897532463, 1.081793370996727, 1.09052883890864, 48)
    6  = @(0.2025961841127114, 0.6844678659950388, 1.26979498242243, 1.126851801446149, 1.120473881934256, 47)
    7  = @(0.1266154959818119, 0.7096750047091113, 1.266381423283964, 1.125336137909009, 1.130546546441805, 46)
    8  = @(0.207816168082475, 0.6898342590982535, 1.300194904219919, 1.140260893050323, 1.133832303844863, 45)
    9  = @(0.1544485185314101, 0.6857633549525511, 1.282012191810668, 1.132259772230148, 1.134644143221363, 44)
    10 = @(0.247847795596123, 0.7276208176551816, 1.343790227660096, 1.159219663247693, 1.145815918573438, 43)
    11 = @(0.181604934442979, 0.7108858875923681, 1.353190151552467, 1.163267016446554, 1.162931656782411, 42)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]  # B
    $ws.Cells.Item($row, 3).Value = $values[1]  # C
    $ws.Cells.Item($row, 4).Value = $values[2]  # D
    $ws.Cells.Item($row, 5).Value = $values[3]  # E
    $ws.Cells.Item($row, 6).Value = $values[4]  # F
    $ws.Cells.Item($row, 7).Value = $values[5]  # G
}
